$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2379.0312
$ws.Range("J17").Value = 2333.4602
$ws.Range("L17").Value = 7000.3806
$ws.Range("N17").Value = -7336.3806
# Row 112
$ws.Range("H112").Value = 1452356.2
$ws.Range("J112").Value = 1590439.1
$ws.Range("L112").Value = 4771317.300000001
$ws.Range("N112").Value = -4773533.300000001
# Row 132
$ws.Range("H132").Value = 1950.7428
$ws.Range("I132").Value = 1893.9395
$ws.Range("K132").Value = 5681.818499999999
$ws.Range("M132").Value = -3151.818499999999
# Row 135
$ws.Range("H135").Value = 1726.037
$ws.Range("I135").Value = 1075.4445
$ws.Range("K135").Value = 9679.0005
$ws.Range("M135").Value = -7144.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2437.9048
$ws.Range("I2").Value = 1966.5834
$ws.Range("J2").Value = 3066.3333
$ws.Range("K2").Value = 1966.5834
$ws.Range("L2").Value = 3066.3333
$ws.Range("M2").Value = -1853.5834
$ws.Range("N2").Value = -3292.3333
# Row 37
$ws.Range("H37").Value = 24034
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
# Row 74
$ws.Range("H74").Value = 2697.762
$ws.Range("I74").Value = 2686.2778
$ws.Range("K74").Value = 2686.2778
$ws.Range("M74").Value = -1812.2778
# Row 77
$ws.Range("H77").Value = 2697.762
$ws.Range("I77").Value = 2686.2778
$ws.Range("K77").Value = 13431.389
$ws.Range("M77").Value = -9063.388999999999
# Row 88
$ws.Range("H88").Value = 1447.5
$ws.Range("I88").Value = 1995
$ws.Range("K88").Value = 1995
$ws.Range("M88").Value = -1589
# Row 91
$ws.Range("H91").Value = 1447.5
$ws.Range("I91").Value = 1995
$ws.Range("K91").Value = 1995
$ws.Range("M91").Value = -591
# Row 116
$ws.Range("H116").Value = 2437.9048
$ws.Range("I116").Value = 1966.5834
$ws.Range("J116").Value = 3066.3333
$ws.Range("K116").Value = 1966.5834
$ws.Range("L116").Value = 3066.3333
$ws.Range("M116").Value = 327.4166
$ws.Range("N116").Value = -7654.3333

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2437.9048
$ws.Range("I3").Value = 1966.5834
$ws.Range("J3").Value = 3066.3333
$ws.Range("K3").Value = 1966.5834
$ws.Range("L3").Value = 3066.3333
$ws.Range("M3").Value = -1852.5834
$ws.Range("N3").Value = -3294.3333
# Row 86
$ws.Range("H86").Value = 2479.6667
$ws.Range("I86").Value = 2763.2727
$ws.Range("K86").Value = 2763.2727
$ws.Range("M86").Value = -1640.2727
# Row 89
$ws.Range("H89").Value = 2479.6667
$ws.Range("I89").Value = 2763.2727
$ws.Range("K89").Value = 13816.3635
$ws.Range("M89").Value = -8200.363499999999
# Row 94
$ws.Range("H94").Value = 1520.8096
$ws.Range("I94").Value = 935.8889
$ws.Range("J94").Value = 1959.5
$ws.Range("K94").Value = 935.8889
$ws.Range("L94").Value = 1959.5
$ws.Range("M94").Value = -484.8889
$ws.Range("N94").Value = -2861.5
# Row 107
$ws.Range("H107").Value = 1422
$ws.Range("I107").Value = 1008.625
$ws.Range("J107").Value = 2248.75
$ws.Range("K107").Value = 1008.625
$ws.Range("L107").Value = 2248.75
$ws.Range("M107").Value = 911.375
$ws.Range("N107").Value = -6088.75
# Row 135
$ws.Range("H135").Value = 92499.75
$ws.Range("J135").Value = 92499.75
$ws.Range("L135").Value = 92499.75
$ws.Range("N135").Value = -102639.75

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 1638.5
$ws.Range("I105").Value = 739.1667
$ws.Range("J105").Value = 2537.8333
$ws.Range("K105").Value = 739.1667
$ws.Range("L105").Value = 2537.8333
$ws.Range("M105").Value = 1007.8333
$ws.Range("N105").Value = -6031.8333
# Row 107
$ws.Range("H107").Value = 1650.85
$ws.Range("I107").Value = 870
$ws.Range("K107").Value = 870
$ws.Range("M107").Value = 1050
# Row 122
$ws.Range("H122").Value = 3172.2092
$ws.Range("I122").Value = 2880.4783
$ws.Range("J122").Value = 3507.7
$ws.Range("K122").Value = 8641.4349
$ws.Range("L122").Value = 10523.1
$ws.Range("M122").Value = -6191.4349
$ws.Range("N122").Value = -15423.1

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 833411.8
$ws.Range("I12").Value = 68.333336
$ws.Range("J12").Value = 1111193
$ws.Range("K12").Value = 205.000008
$ws.Range("L12").Value = 3333579
$ws.Range("M12").Value = -32.00000800000001
$ws.Range("N12").Value = -3333925
# Row 93
$ws.Range("H93").Value = 1224.6666
$ws.Range("I93").Value = 449.33334
$ws.Range("K93").Value = 1348.00002
$ws.Range("M93").Value = 523.9999800000001

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1202.25
$ws.Range("I97").Value = 802.5714
$ws.Range("K97").Value = 802.5714
$ws.Range("M97").Value = -306.5714
# Row 102
$ws.Range("H102").Value = 3359.8
$ws.Range("I102").Value = 3199.75
$ws.Range("K102").Value = 3199.75
$ws.Range("M102").Value = -1577.75
# Row 107
$ws.Range("H107").Value = 3083.1667
$ws.Range("I107").Value = 2499
$ws.Range("J107").Value = 3200
$ws.Range("K107").Value = 2499
$ws.Range("L107").Value = 3200
$ws.Range("M107").Value = -579
$ws.Range("N107").Value = -7040
# Row 122
$ws.Range("H122").Value = 1544.8846
$ws.Range("I122").Value = 1507.8125
$ws.Range("J122").Value = 1604.2
$ws.Range("K122").Value = 4523.4375
$ws.Range("L122").Value = 4812.6
$ws.Range("M122").Value = -2073.4375
$ws.Range("N122").Value = -9712.6
# Row 132
$ws.Range("H132").Value = 2616.9119
$ws.Range("I132").Value = 2316.4482
$ws.Range("J132").Value = 4359.6
$ws.Range("K132").Value = 6949.344599999999
$ws.Range("L132").Value = 13078.8
$ws.Range("M132").Value = -4419.344599999999
$ws.Range("N132").Value = -18138.8

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6178.7856
$ws.Range("I7").Value = 6046.1816
$ws.Range("K7").Value = 6046.1816
$ws.Range("M7").Value = -5934.1816
# Row 22
$ws.Range("H22").Value = 3439.2222
$ws.Range("I22").Value = 2617.1667
$ws.Range("K22").Value = 2617.1667
$ws.Range("M22").Value = -2322.1667
# Row 27
$ws.Range("H27").Value = 3439.2222
$ws.Range("I27").Value = 2617.1667
$ws.Range("K27").Value = 2617.1667
$ws.Range("M27").Value = -2510.1667
# Row 40
$ws.Range("H40").Value = 4611.5
$ws.Range("I40").Value = 3950.5386
$ws.Range("J40").Value = 6330
$ws.Range("K40").Value = 3950.5386
$ws.Range("L40").Value = 6330
$ws.Range("M40").Value = -3814.5386
$ws.Range("N40").Value = -6602
# Row 46
$ws.Range("H46").Value = 7217.5938
$ws.Range("J46").Value = 7745.107
$ws.Range("L46").Value = 7745.107
$ws.Range("N46").Value = -8121.107
# Row 61
$ws.Range("H61").Value = 1866.5
$ws.Range("I61").Value = 1344.3334
$ws.Range("K61").Value = 1344.3334
$ws.Range("M61").Value = -1142.3334
# Row 82
$ws.Range("H82").Value = 2498.4167
$ws.Range("I82").Value = 2736.5334
$ws.Range("J82").Value = 2101.5557
$ws.Range("K82").Value = 2736.5334
$ws.Range("L82").Value = 2101.5557
$ws.Range("M82").Value = -2375.5334
$ws.Range("N82").Value = -2823.5557
# Row 85
$ws.Range("H85").Value = 2498.4167
$ws.Range("I85").Value = 2736.5334
$ws.Range("J85").Value = 2101.5557
$ws.Range("K85").Value = 2736.5334
$ws.Range("L85").Value = 2101.5557
$ws.Range("M85").Value = -1488.5334
$ws.Range("N85").Value = -4597.5557
# Row 93
$ws.Range("H93").Value = 2587.8948
$ws.Range("I93").Value = 2075.1538
$ws.Range("J93").Value = 3698.8333
$ws.Range("K93").Value = 2075.1538
$ws.Range("L93").Value = 3698.8333
$ws.Range("M93").Value = -827.1538
$ws.Range("N93").Value = -6194.8333
# Row 100
$ws.Range("H100").Value = 1748.5
$ws.Range("I100").Value = 1495
$ws.Range("K100").Value = 1495
$ws.Range("M100").Value = -954
# Row 113
$ws.Range("H113").Value = 1866.5
$ws.Range("I113").Value = 1344.3334
$ws.Range("K113").Value = 1344.3334
$ws.Range("M113").Value = 825.6666
# Row 126
$ws.Range("H126").Value = 6178.7856
$ws.Range("I126").Value = 6046.1816
$ws.Range("K126").Value = 18138.5448
$ws.Range("M126").Value = -15668.5448

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 7687.75
$ws.Range("I122").Value = 8500.333000000001
$ws.Range("K122").Value = 25500.999
$ws.Range("M122").Value = -23050.999

Write-Host "Applied market-price refresh to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
